$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.693.90'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +3.52%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.423.29'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.11%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.56'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.78%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.70'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +6.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.514'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.25%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.527'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +9.26%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.53'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.25%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.67%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.92'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.87%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.46%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.03%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.801.45'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.35%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.426.47'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.67%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +4.24%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '44.543.89'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.23%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.31'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.38'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.25%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.77'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.84%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '242.19'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.75%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.03%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.70%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.10%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.20'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.47%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.44%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.63%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.60'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.37%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '48.56'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.67%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.126'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +18.46%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.55'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +11.00%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0780'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +8.57%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.14%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.28%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.45%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.72%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.68%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '121.72'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.48%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.110'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.76%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.99%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.04'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.08%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0290'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.97%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.946.40'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.42%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.34%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +8.40%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.90%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.70'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +11.87%  '

$ws.Range("B50").Value = 'BitcoinSV'
$ws.Range("C50").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '75.23'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.81%  '

$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '54.27'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +5.65%  '
